$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph whose text starts with a given needle so the
# script does not depend on fixed paragraph indices that could shift.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($needle) {
    $idx = 1
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($needle)) {
            return $idx
        }
        $idx = $idx + 1
    }
    throw "Find-ParagraphIndex: no paragraph found starting with '" + $needle + "'"
}

# ===========================================================================
# Edit 1: insert two new bullet paragraphs (ilvl=2) right before the
# paragraph that begins with "Το σύστημα εμφανίζει το κόστος".
# ===========================================================================
$costParaIdx = Find-ParagraphIndex("Το σύστημα εμφανίζει το κόστος")
$costPara = $d.Paragraphs.Item($costParaIdx)

# Insert two empty paragraphs immediately before it (they inherit the
# surrounding list formatting/level/language automatically).
$costPara.Range.InsertParagraphBefore()
$costPara.Range.InsertParagraphBefore()

$newPara1 = $d.Paragraphs.Item($costParaIdx)
$newPara2 = $d.Paragraphs.Item($costParaIdx + 1)

$newPara1.Range.Text = "Το σύστημα εμφανίζει λίστα με αντικείμενα που πουλάνε άλλοι παίκτες."
$newPara2.Range.Text = "Ο παίκτης επιλέγει το αντικείμενο που επιθημεί να αγοράσει."

# ===========================================================================
# Edit 2: insert a new bullet paragraph (ilvl=3) right after the paragraph
# that begins with "Ο παίκτης επιλέγει την τιμή της προσφοράς του." (and
# before "Εντός 24 ωρών ...").
# ===========================================================================
$priceParaIdx = Find-ParagraphIndex(" Ο παίκτης επιλέγει την τιμή της προσφοράς")
$pricePara = $d.Paragraphs.Item($priceParaIdx)
$pricePara.Range.InsertParagraphAfter()

$newPara3Idx = $priceParaIdx + 1
$newPara3 = $d.Paragraphs.Item($newPara3Idx)
$newPara3.Range.Text = "Το σύστημα ανανεώνει την νέα μέγιστη προσφορά στον "

# Grab a zero-formatting (no w:lang, plain Latin) run from elsewhere in the
# document -- e.g. the "shop" run in paragraph 1 -- and splice a COPY of its
# formatting in at the end of our new paragraph, then rename the copied text
# to "server". This produces a separate run with no rPr/lang, matching how
# plain Latin runs look throughout this document.
$shopSearch = $d.Content
$shopSearch.Find.Execute("shop", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$shopFormatted = $shopSearch.FormattedText

$beforeLen = $newPara3.Range.End
$insertPoint = $d.Range($newPara3.Range.End - 1, $newPara3.Range.End - 1)
$insertPoint.FormattedText = $shopFormatted
$afterLen = $newPara3.Range.End

$serverRange = $d.Range($beforeLen - 1, $afterLen - 1)
$serverRange.Text = "server"

# Grab a lone Greek-language run (no extra rsid attribute) from paragraph 6
# ("... Η ροή συνεχίζεται στο βήμα 6.  ") and splice a copy of ITS formatting
# in at the very end of the new paragraph, then rename the copied text to
# ".". This produces a third run carrying w:lang="el-GR" again.
$trailingSpaceParaIdx = Find-ParagraphIndex("Ο παίκτης πληρώνει τ")
$trailingSpacePara = $d.Paragraphs.Item($trailingSpaceParaIdx)
$spaceRunRange = $d.Range($trailingSpacePara.Range.End - 2, $trailingSpacePara.Range.End - 1)
$spaceFormatted = $spaceRunRange.FormattedText

$dotInsertPoint = $d.Range($newPara3.Range.End - 1, $newPara3.Range.End - 1)
$dotInsertPoint.FormattedText = $spaceFormatted
$dotFinalEnd = $newPara3.Range.End

$dotRange = $d.Range($dotFinalEnd - 2, $dotFinalEnd - 1)
$dotRange.Text = "."

Write-Host "Edit complete."
